$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 294 (shifts old row 294 and below down by one)
# New row: Durazno, Florida King, Especial  (Feria Lagunitas de Puerto Montt, 2020-12-11)
$ws.Rows.Item(294).Insert()

$ws.Cells.Item(294, 1).Value = 4
$ws.Cells.Item(294, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(294, 3).Value = 'Los Lagos'
$ws.Cells.Item(294, 4).Value = 44176
$ws.Cells.Item(294, 5).Value = 10
$ws.Cells.Item(294, 6).Value = 'Fruta'
$ws.Cells.Item(294, 7).Value = 100103
$ws.Cells.Item(294, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(294, 9).Value = 100103004
$ws.Cells.Item(294, 10).Value = 'Durazno'
$ws.Cells.Item(294, 11).Value = 'Florida King'
$ws.Cells.Item(294, 12).Value = 'Especial'
$ws.Cells.Item(294, 13).Value = 300
$ws.Cells.Item(294, 14).Value = 24000
$ws.Cells.Item(294, 15).Value = 24000
$ws.Cells.Item(294, 16).Value = 24000
$ws.Cells.Item(294, 17).Value = '$/caja 12 kilos empedrada'
$ws.Cells.Item(294, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(294, 19).Value = 2000
$ws.Cells.Item(294, 20).Value = 12

# Insert a new row at 297 (shifts old Springcrest-Primera and below down by one)
# New row: Durazno, Springcrest, Especial  (Feria Lagunitas de Puerto Montt, 2020-12-11)
$ws.Rows.Item(297).Insert()

$ws.Cells.Item(297, 1).Value = 4
$ws.Cells.Item(297, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(297, 3).Value = 'Los Lagos'
$ws.Cells.Item(297, 4).Value = 44176
$ws.Cells.Item(297, 5).Value = 10
$ws.Cells.Item(297, 6).Value = 'Fruta'
$ws.Cells.Item(297, 7).Value = 100103
$ws.Cells.Item(297, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(297, 9).Value = 100103004
$ws.Cells.Item(297, 10).Value = 'Durazno'
$ws.Cells.Item(297, 11).Value = 'Springcrest'
$ws.Cells.Item(297, 12).Value = 'Especial'
$ws.Cells.Item(297, 13).Value = 300
$ws.Cells.Item(297, 14).Value = 24000
$ws.Cells.Item(297, 15).Value = 24000
$ws.Cells.Item(297, 16).Value = 24000
$ws.Cells.Item(297, 17).Value = '$/caja 12 kilos empedrada'
$ws.Cells.Item(297, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(297, 19).Value = 2000
$ws.Cells.Item(297, 20).Value = 12

Write-Output 'done'
